$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.065860333333333
$ws.Range("H2").Value = 6.197581
$ws.Range("I2").Value = 0.546612493277129
$ws.Range("J2").Value = 0.546612493277129
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.065860333333333
$ws.Range("N2").Value = 6.197581
$ws.Range("O2").Value = 0.546612493277129
$ws.Range("P2").Value = 0.546612493277129
$ws.Range("Q2").Value = 4.26777891684011
$ws.Range("R2").Value = 38.410010251561
$ws.Range("S2").Value = 0.2987852178066394
$ws.Range("T2").Value = 0.2987852178066394

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.065860333333333
$ws.Range("H3").Value = 6.197581
$ws.Range("I3").Value = 0.546612493277129
$ws.Range("J3").Value = 0.546612493277129
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.8953543333333333
$ws.Range("N3").Value = 2.686063
$ws.Range("O3").Value = 0.2369046235183445
$ws.Range("P3").Value = 0.2369046235183445
$ws.Range("Q3").Value = 1.849677001511444
$ws.Range("R3").Value = 16.647093013603
$ws.Range("S3").Value = 0.1294950269302419
$ws.Range("T3").Value = 0.1294950269302419

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.065860333333333
$ws.Range("H4").Value = 6.197581
$ws.Range("I4").Value = 0.546612493277129
$ws.Range("J4").Value = 0.546612493277129
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.8181726666666668
$ws.Range("N4").Value = 2.454518
$ws.Range("O4").Value = 0.2164828832045265
$ws.Range("P4").Value = 0.2164828832045265
$ws.Range("Q4").Value = 1.690230457884222
$ws.Range("R4").Value = 15.212074120958
$ws.Range("S4").Value = 0.1183322485402478
$ws.Range("T4").Value = 0.1183322485402477

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.8953543333333333
$ws.Range("H5").Value = 2.686063
$ws.Range("I5").Value = 0.2369046235183445
$ws.Range("J5").Value = 0.2369046235183445
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.065860333333333
$ws.Range("N5").Value = 6.197581
$ws.Range("O5").Value = 0.546612493277129
$ws.Range("P5").Value = 0.546612493277129
$ws.Range("Q5").Value = 1.849677001511444
$ws.Range("R5").Value = 16.647093013603
$ws.Range("S5").Value = 0.1294950269302419
$ws.Range("T5").Value = 0.1294950269302419

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.8953543333333333
$ws.Range("H6").Value = 2.686063
$ws.Range("I6").Value = 0.2369046235183445
$ws.Range("J6").Value = 0.2369046235183445
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.8953543333333333
$ws.Range("N6").Value = 2.686063
$ws.Range("O6").Value = 0.2369046235183445
$ws.Range("P6").Value = 0.2369046235183445
$ws.Range("Q6").Value = 0.8016593822187776
$ws.Range("R6").Value = 7.214934439968999
$ws.Range("S6").Value = 0.05612380064436854
$ws.Range("T6").Value = 0.05612380064436854

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.8953543333333333
$ws.Range("H7").Value = 2.686063
$ws.Range("I7").Value = 0.2369046235183445
$ws.Range("J7").Value = 0.2369046235183445
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.8181726666666668
$ws.Range("N7").Value = 2.454518
$ws.Range("O7").Value = 0.2164828832045265
$ws.Range("P7").Value = 0.2164828832045265
$ws.Range("Q7").Value = 0.7325544425148889
$ws.Range("R7").Value = 6.592989982634
$ws.Range("S7").Value = 0.05128579594373409
$ws.Range("T7").Value = 0.05128579594373409

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.8181726666666668
$ws.Range("H8").Value = 2.454518
$ws.Range("I8").Value = 0.2164828832045265
$ws.Range("J8").Value = 0.2164828832045265
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.065860333333333
$ws.Range("N8").Value = 6.197581
$ws.Range("O8").Value = 0.546612493277129
$ws.Range("P8").Value = 0.546612493277129
$ws.Range("Q8").Value = 1.690230457884222
$ws.Range("R8").Value = 15.212074120958
$ws.Range("S8").Value = 0.1183322485402478
$ws.Range("T8").Value = 0.1183322485402477

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.8181726666666668
$ws.Range("H9").Value = 2.454518
$ws.Range("I9").Value = 0.2164828832045265
$ws.Range("J9").Value = 0.2164828832045265
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.8953543333333333
$ws.Range("N9").Value = 2.686063
$ws.Range("O9").Value = 0.2369046235183445
$ws.Range("P9").Value = 0.2369046235183445
$ws.Range("Q9").Value = 0.7325544425148889
$ws.Range("R9").Value = 6.592989982634
$ws.Range("S9").Value = 0.05128579594373409
$ws.Range("T9").Value = 0.05128579594373409

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.8181726666666668
$ws.Range("H10").Value = 2.454518
$ws.Range("I10").Value = 0.2164828832045265
$ws.Range("J10").Value = 0.2164828832045265
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.8181726666666668
$ws.Range("N10").Value = 2.454518
$ws.Range("O10").Value = 0.2164828832045265
$ws.Range("P10").Value = 0.2164828832045265
$ws.Range("Q10").Value = 0.6694065124804446
$ws.Range("R10").Value = 6.024658612324001
$ws.Range("S10").Value = 0.04686483872054466
$ws.Range("T10").Value = 0.04686483872054464
